$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elem = $wb.Worksheets.Item("Elements")

$newUrl = "https://johnmoehrke.github.io/MHV-PGHD/StructureDefinition/FM-Race"
$newName = "Race"
$newDate = "2022-04-11T07:37:02-05:00"
$newDesc = "What is the Race of this family member`n`nNote would like to use us-core defined extension, but is not allowed in FamilyMemberHistory [Jira issue](https://jira.hl7.org/browse/FHIR-35998)"
$newValueSet = "http://terminology.hl7.org/ValueSet/v3-Race"

# Metadata sheet
$meta.Range("B2").Value = $newUrl
$meta.Range("B4").Value = $newName
$meta.Range("B5").Value = $newName
$meta.Range("B8").Value = $newDate
$meta.Range("B12").Value = $newDesc

# Elements sheet
$elem.Range("K2").Value = $newName
$elem.Range("L2").Value = $newDesc
$elem.Range("Q5").Value = $newUrl
$elem.Range("Y7").Value = $newValueSet

# Refresh the "Binding Value Set" column (Y, column 25) width to fit the
# shorter value-set URL text now shown in Y7 (narrows from ~51.7 to ~41.5 chars)
$elem.Range("Y1").ColumnWidth = 40.666666
